# Harmonize tag TANs between templates
#
# Target changes (per commit "harmonize tag TANs between templates"):
#  - isa_template!E13: "sampling" -> "Sampling"
#  - isa_template!E14: "http://purl.obolibrary.org/obo/C_6774"
#                      -> "http://purl.obolibrary.org/obo/NCIT_C25662"
#  - isa_template!E15: cleared (was "C")
#  - isa_template!E17: cleared (was "false")
#  - selection / scroll position cosmetic updates on isa_template sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")
$ws.Activate()

# Tags row: capitalize "Sampling"
$ws.Range("E13").Value = "Sampling"

# Tags Term Accession Number row: swap in the harmonized NCIT accession
$ws.Range("E14").Value = "http://purl.obolibrary.org/obo/NCIT_C25662"

# Tags Term Source REF row: drop the now-unused 5th tag's source ref
$ws.Range("E15").ClearContents()

# Comment[isObsolete] row: drop the now-unused 5th tag's obsolete flag
$ws.Range("E17").ClearContents()

# Reflect the author's final cursor/scroll position on save
$ws.Range("F13").Select()
